$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark the "Snow" feature row as implemented: D13 0 -> 10, and label it
# "DONE" in E13 (matching the styling already used on the other
# completed-feature rows, e.g. E4/E5/E8/E9/E14/E16).
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = "DONE"
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the view/selection that was active when the sheet was saved.
$ws.Activate()
$ws.Range("E11").Select()
